# Plantilla de Casos de Uso - iteracion 2 update
# - Corrige columna "Incremento" (G): de "Documento ERS" a "Primer parcial" para todos los CU
# - Corrige columna "Esfuerzo (hrs)" (F): de "01:00 hr" a nuevos valores por CU
# - Actualiza la celda seleccionada / vista de la hoja

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# Nuevos valores de esfuerzo por fila (columna F) e incremento (columna G)
$effort = @{
    5  = "05:30 hrs"
    6  = "04:30 hrs"
    7  = "03:30 hrs"
    8  = "03:30 hrs"
    9  = "05:30 hrs"
    10 = "04:00 hrs"
    11 = "03:30 hr"
    12 = "03:30 hrs"
    13 = "03:30 hrs"
    14 = "04:00 hrs"
    15 = "03:30 hr"
    16 = "03:30 hrs"
    17 = "03:30 hrs"
    18 = "04:30 hrs"
    19 = "03:30 hrs"
    20 = "03:00 hrs"
    21 = "04:30 hrs"
    22 = "03:30 hrs"
    23 = "03:00 hrs"
}

for ($row = 5; $row -le 23; $row++) {
    $ws.Cells.Item($row, 7).Value = "Primer parcial"
    $ws.Cells.Item($row, 6).Value = $effort[$row]
}

# Actualiza la selección de la hoja activa
$ws.Range("D19").Select()
